# "add work dowload button"
# The tar_nodes graph-data sheet drops the two stale/duplicate
# "...level11.xlsx" stub nodes (ids 12 and 13, rows 4 and 8) and wires up
# real stats for the `download_files` work node (id 3) which previously had
# no seconds/bytes/branch/group data - i.e. gives it a working download
# button. A few sibling nodes' seconds/bytes figures are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete stub rows. Row 4 (id=12) first; after that row 8
# (id=13) has already shifted up to row 7, so delete it there.
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(7).Delete() | Out-Null

# Refresh seconds/bytes on the correlation_graph (id=1) node.
$ws.Range("E2").Value = 0.613
$ws.Range("F2").Value = 163561

# Populate the download_files (id=3) node with real stats/group info -
# this is the "work download button" becoming active.
$ws.Range("E4").Value = 0.003
$ws.Range("F4").Value = 13922
$ws.Range("J4").Value = -289
$ws.Range("K4").Value = -30
$ws.Range("L4").Value = "download"
$ws.Range("M4").Value = "Загрузка исходных данных"
$ws.Range("N4").Value = "Блок загрузки"

# Refresh seconds/bytes on the model_graph (id=8) node.
$ws.Range("E10").Value = 0.319
$ws.Range("F10").Value = 225872

# Refresh seconds on the send_report (id=10) node.
$ws.Range("E12").Value = 0.001
